$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$BVals = @(1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02)
$CVals = @(1.018970544301157, 1.020058099956599, 1.020761729821542, 1.021057514453097, 1.021107176756319, 1.020765682193571, 1.019338108052018, 1.01682181693209, 1.015143767717658, 1.014417021235507, 1.014147053306122, 1.014204963349799, 1.01439470605465, 1.014511609865488, 1.015191996461443, 1.015618746869038, 1.015867649703282, 1.015952516881473, 1.01557296198726, 1.014338832169656, 1.013562758469114, 1.013974182182249, 1.015593650245808, 1.017472427855695)
$DVals = @(1.024472714553104, 1.025264769683957, 1.025776925117156, 1.025992149204427, 1.026028281213447, 1.025779801290781, 1.024740466928942, 1.022906314639513, 1.02168175192038, 1.021151082468404, 1.02095390458632, 1.02099620278489, 1.021134784979037, 1.021220161622491, 1.021716961763951, 1.022028477632931, 1.022210138471221, 1.022272073141273, 1.021995059176707, 1.021093977730443, 1.020527064142488, 1.020827630524416, 1.022010159665498, 1.023380805378414)
$EVals = @(1.020160889744303, 1.021086997280732, 1.021686551765207, 1.021938676251165, 1.021981013262987, 1.02168992038323, 1.020473809637924, 1.018333183833983, 1.016907670418715, 1.016290781685331, 1.016061697074399, 1.016110833996509, 1.016271844350923, 1.016371055477601, 1.016948619049281, 1.017311007972163, 1.017522418801804, 1.017594510451663, 1.017272123397749, 1.016224429286845, 1.015566022568765, 1.015915026066182, 1.017289693553336, 1.018886311036184)
$IVals = @(1.02735632998126, 1.027515291527106, 1.027616625202157, 1.027658860517189, 1.027665930579961, 1.02761719098845, 1.027410367403291, 1.02703424865245, 1.026775671120321, 1.026661849248019, 1.026619292037362, 1.026628433312723, 1.026658337140184, 1.026676724949566, 1.026783186010844, 1.026849469404599, 1.026887952290021, 1.026901043576732, 1.02684237634738, 1.026649538900535, 1.026526681952746, 1.026591963528981, 1.026845581947964, 1.027132865774103)
$JVals = @(1.024175959226467, 1.024899344706809, 1.025366817198803, 1.025563197383106, 1.025596161959458, 1.025369441811315, 1.024420556209458, 1.022743859192415, 1.02162294477139, 1.02113683701112, 1.020956162974871, 1.020994923237954, 1.021121904723838, 1.021200127367158, 1.02165519044743, 1.021940440102604, 1.022106749583637, 1.022163444641373, 1.021909842943339, 1.021084514933153, 1.02056495056518, 1.020840442923761, 1.021923668704317, 1.023177875043669)
$KVals = @(1.027301114536878, 1.027900309103592, 1.028287020838992, 1.028449352993378, 1.028476595107474, 1.028289190878378, 1.02750382350003, 1.026112205784765, 1.025179300241799, 1.024774120654594, 1.024623434953166, 1.024655765870261, 1.02476166867954, 1.024826894522633, 1.02520616487144, 1.025443743178296, 1.025582200350001, 1.025629390608755, 1.025418265502963, 1.024730488018994, 1.024296990730721, 1.024526896684915, 1.025429778132994, 1.026472883082035)
$LVals = @(1.023001995971831, 1.023733954080661, 1.024207260901449, 1.024406162546801, 1.024439554514879, 1.024209918935945, 1.023249430597007, 1.021554490981639, 1.0204228998822, 1.019932522126689, 1.019750315151428, 1.019789401829701, 1.019917462041937, 1.019996356366491, 1.020455436389648, 1.020743300353251, 1.020911168777509, 1.020968401169511, 1.020712419179728, 1.019879753126, 1.019355882457509, 1.019633628500526, 1.020726373169303, 1.021992961519722)
$NVals = @(1.012108358176944, 1.012351135332999, 1.01250792359418, 1.012573764309795, 1.012584814967867, 1.012508803648287, 1.01219046890246, 1.011627192318421, 1.011250112818804, 1.011086463528539, 1.011025621121942, 1.011038674548237, 1.011081435413876, 1.011107774404551, 1.011260965860602, 1.011356959401016, 1.011412914975044, 1.011431988301011, 1.011346663908886, 1.011068844944192, 1.010893846454701, 1.010986647044536, 1.011351316109295, 1.011773088438938)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value2 = $BVals[$i]
    $ws.Cells.Item($row, 3).Value2 = $CVals[$i]
    $ws.Cells.Item($row, 4).Value2 = $DVals[$i]
    $ws.Cells.Item($row, 5).Value2 = $EVals[$i]
    $ws.Cells.Item($row, 9).Value2 = $IVals[$i]
    $ws.Cells.Item($row, 10).Value2 = $JVals[$i]
    $ws.Cells.Item($row, 11).Value2 = $KVals[$i]
    $ws.Cells.Item($row, 12).Value2 = $LVals[$i]
    $ws.Cells.Item($row, 14).Value2 = $NVals[$i]
}
